$wb = $excel.ActiveWorkbook

# --- BATT_DCDC sheet: update Package for positions 21, 27, 28 and 35 from SMD 0402 to SMD 0603 ---
$wsBom = $wb.Worksheets.Item("BATT_DCDC")

$wsBom.Range("G27").Value = "SMD 0603"
$wsBom.Range("G33").Value = "SMD 0603"
$wsBom.Range("G34").Value = "SMD 0603"
$wsBom.Range("G41").Value = "SMD 0603"

# --- Clear the "NEW" flag for positions that are no longer new (position 30 and 38) ---
$wsBom.Range("B36").Value = "NO"
$wsBom.Range("B44").Value = "NO"

# Position 30 used to be manually highlighted (yellow fill) while it was NEW;
# drop that highlight now that it's no longer new, matching the other rows.
$wsBom.Range("B35").Copy()
$wsBom.Range("B36").PasteSpecial(-4122) # xlPasteFormats

# --- _HISTORY sheet: log this change as a new version ---
$wsHistory = $wb.Worksheets.Item("_HISTORY")

$wsHistory.Range("A9").Value = 6
$wsHistory.Range("B9").Value = (Get-Date -Year 2023 -Month 5 -Day 16).Date
$wsHistory.Range("C9").Value = "JRC"
$wsHistory.Range("D9").Value = "Las posiciones 21,27,28 y 35 cambian su Package a SMD 0603"

# Match formatting used by the rest of the history table (date format + wrap text)
$wsHistory.Range("B9").NumberFormat = $wsHistory.Range("B8").NumberFormat
$wsHistory.Range("D4:D9").WrapText = $true
$wsHistory.Range("D3").WrapText = $true
$wsHistory.Range("D7").RowHeight = 30
$wsHistory.Range("D9").RowHeight = 30

$wb.Save()
